$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need NumberFormat forced
# to Text ("@") before assignment, otherwise Excel auto-converts the string to a
# number (e.g. "8.50" -> 8.5, losing the trailing zero / exact textual formatting).

$ws.Range('D2').Value = '57.955.98'
$ws.Range('E2').Value = '  -0.54%  '

$ws.Range('D3').Value = '2.357.92'
$ws.Range('E3').Value = '  -0.50%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.28'
$ws.Range('E5').Value = '  -0.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.62'
$ws.Range('E6').Value = '  -0.48%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  +4.44%  '

$ws.Range('E9').Value = '  -0.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.54'
$ws.Range('E10').Value = '  +1.93%  '

$ws.Range('E11').Value = '  -1.35%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  +0.07%  '

$ws.Range('D13').Value = '2.778.30'
$ws.Range('E13').Value = '  -0.18%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.77'
$ws.Range('E14').Value = '  +0.16%  '

$ws.Range('D15').Value = '57.953.86'
$ws.Range('E15').Value = '  -0.36%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000134'
$ws.Range('E16').Value = '  +0.23%  '

$ws.Range('D17').Value = '2.374.79'
$ws.Range('E17').Value = '  +0.71%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.79'
$ws.Range('E18').Value = '  +2.40%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '330.56'
$ws.Range('E19').Value = '  -2.71%  '

$ws.Range('E20').Value = '  +1.58%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.71'
$ws.Range('E21').Value = '  -2.64%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '62.75'
$ws.Range('E23').Value = '  +0.75%  '

$ws.Range('E24').Value = '  -1.61%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.50'
$ws.Range('E25').Value = '  -0.90%  '

$ws.Range('E26').Value = '  +0.91%  '

$ws.Range('E27').Value = '  -2.66%  '

$ws.Range('E28').Value = '  -0.19%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.71'
$ws.Range('E29').Value = '  -2.50%  '

$ws.Range('D30').Value = '0.0₃0741'
$ws.Range('E30').Value = '  +0.39%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.14'
$ws.Range('E31').Value = '  -0.27%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.04'

$ws.Range('E33').Value = '  -1.01%  '

$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.23'
$ws.Range('E35').Value = '  +3.00%  '

$ws.Range('E36').Value = '  +0.06%  '

$ws.Range('E37').Value = '  -2.31%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.61'
$ws.Range('E38').Value = '  +0.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.43'
$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '143.05'
$ws.Range('E40').Value = '  -4.56%  '

$ws.Range('E41').Value = '  +0.49%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.66'
$ws.Range('E42').Value = '  +0.34%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '289.62'
$ws.Range('E43').Value = '  +1.52%  '

$ws.Range('E44').Value = '  +1.60%  '

$ws.Range('E45').Value = '  +0.74%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.16'
$ws.Range('E46').Value = '  +0.93%  '

$ws.Range('E47').Value = '  +0.69%  '

$ws.Range('E48').Value = '  +2.07%  '

$ws.Range('B49').Value = 'Polygon'
$ws.Range('C49').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.382'
$ws.Range('E49').Value = '  -0.28%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.50'
$ws.Range('E50').Value = '  -0.57%  '

$ws.Range('E51').Value = '  +1.61%  '
